$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 23: add dating method "AMS" to column B (copy format from B22) ---
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B23").Value = "AMS"

# --- Row 24: new site "Hajinri" ---
$ws.Range("A24").Value = "Hajinri"

$ws.Range("B22").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = "AMS"

$ws.Range("C22").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "42,000±340"

$ws.Range("D22").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = 127.8

$ws.Range("G22").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("G24").Value = 86

$ws.Range("H22").Copy()
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("H24").Value = "한국선사문화연구원 2018"

$ws.Rows.Item(24).RowHeight = 15.75

$ws.Application.CutCopyMode = $false

# --- Update active selection to match latest edit location ---
$ws.Range("H25").Select()
